$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from row 2 to row 289: increment date serial 45178 -> 45179
$ws.Range("C2:C289").Value = 45179
